# Update vm_pu results: case with 380 kV done.
# Voltage setpoint (column B) changed from 1.05 to 1.02, and the recalculated
# bus voltage magnitudes (columns C:F, I:N) for data rows 2-25 are updated to match.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @{
    2 = @{2=1.02;3=1.072714666073793;4=1.075858415397501;5=1.085367558967875;6=1.090296748170131;9=1.054828572189578;10=1.077632818300618;11=1.078543875911341;12=1.088028171411269;13=1.092944664306481;14=1.029224888212585}
    3 = @{2=1.02;3=1.07419035827153;4=1.077044819703253;5=1.086776655987513;6=1.09165501277119;9=1.055253965002003;10=1.078764458076185;11=1.079546771329144;12=1.089255028720666;13=1.09412173580405;14=1.029626136867439}
    4 = @{2=1.02;3=1.07514413091314;4=1.077811341319072;5=1.087687635687663;6=1.092533016948944;9=1.05552713379347;10=1.079495145909311;11=1.080193973291176;12=1.090047551146669;13=1.094881948317036;14=1.029884697867383}
    5 = @{2=1.02;3=1.075544840179138;4=1.078133313246248;5=1.088070425270012;6=1.092901923272243;9=1.055641475880192;10=1.07980195857496;11=1.080465644823173;12=1.090380412529164;13=1.095201203544755;14=1.02999314147903}
    6 = @{2=1.02;3=1.075612106155291;4=1.078187357819455;5=1.088134686555313;6=1.09296385229385;9=1.055660645270038;10=1.079853452298003;11=1.080511235614983;12=1.090436283086498;13=1.095254788207865;14=1.030011334706745}
    7 = @{2=1.02;3=1.075149486210974;4=1.077815644594136;5=1.087692751266175;6=1.092537947099064;9=1.055528663592662;10=1.079499246997322;11=1.080197604994076;12=1.090052000088634;13=1.094886215543892;14=1.029886147897894}
    8 = @{2=1.02;3=1.073213614506562;4=1.076259609312599;5=1.085843938811924;6=1.090755966050358;9=1.054972769384661;10=1.078015586834779;11=1.078883171489014;12=1.088443073206152;13=1.093342760285983;14=1.02936071536741}
    9 = @{2=1.02;3=1.069793661216129;4=1.073508590551309;5=1.082579727219062;6=1.087608895548132;9=1.0539771253168;10=1.075389039710422;11=1.076553476889224;12=1.085597485431732;13=1.090611825738201;14=1.028426540177754}
    10 = @{2=1.02;3=1.067507444382811;4=1.071668199001212;5=1.080398969152825;6=1.085505833225579;9=1.053302419416842;10=1.073629544886519;11=1.074991019533983;12=1.08369309113587;13=1.088783403051033;14=1.027798076665373}
    11 = @{2=1.02;3=1.066515916117314;4=1.070869715255135;5=1.079453506419089;6=1.084593923855108;9=1.053007640039621;10=1.072865590458195;11=1.074312188595502;12=1.082866658337639;13=1.087989760649999;14=1.027524574760148}
    12 = @{2=1.02;3=1.066147373039342;4=1.070572879816084;5=1.079102135592657;6=1.084255003211133;9=1.052897748592029;10=1.072581505238174;11=1.074059693503179;12=1.082559405301284;13=1.08769467186014;14=1.027422775660704}
    13 = @{2=1.02;3=1.066226438007566;4=1.070636563108814;5=1.079177514254229;6=1.084327711849506;9=1.052921338698564;10=1.072642456989915;11=1.074113870343708;12=1.08262532484776;13=1.087757982833584;14=1.027444621372447}
    14 = @{2=1.02;3=1.066485457281297;4=1.070845183733756;5=1.079424465786049;6=1.084565912621526;9=1.05299856449709;10=1.072842114416892;11=1.074291324381451;12=1.082841266434403;13=1.087965374562823;14=1.02751616427615}
    15 = @{2=1.02;3=1.066645014862623;4=1.070973689422877;5=1.079576596213218;6=1.084712649702122;9=1.053046093142844;10=1.072965087564678;11=1.074400613531509;12=1.082974278008336;13=1.088093116234529;14=1.027560216555818}
    16 = @{2=1.02;3=1.067573214870942;4=1.071721157950892;5=1.080461690977206;6=1.085566326374317;9=1.053321927382994;10=1.073680201675701;11=1.07503602291445;12=1.083747899951687;13=1.088836033475006;14=1.027816198984519}
    17 = @{2=1.02;3=1.068155021392616;4=1.072189597892534;5=1.081016566964843;6=1.086101470969652;9=1.053494245729286;10=1.074128212489319;11=1.075433985388195;12=1.084232681710011;13=1.089301526888522;14=1.027976400964694}
    18 = @{2=1.02;3=1.06849422724919;4=1.072462678702554;5=1.081340103325192;6=1.086413489689128;9=1.053594502785097;10=1.074389329171315;11=1.075665890849856;12=1.084515271756023;13=1.089572856137846;14=1.028069711727134}
    19 = @{2=1.02;3=1.06860986212486;4=1.072555766532483;5=1.081450401823775;6=1.086519859489873;9=1.053628644959252;10=1.074478329387506;11=1.075744927570629;12=1.084611598165855;13=1.089665341161601;14=1.02810150589926}
    20 = @{2=1.02;3=1.068092614841969;4=1.072139354510135;5=1.080957045801498;6=1.086044067690135;9=1.053475783822111;10=1.074080165944116;11=1.075391310502443;12=1.084180687355208;13=1.089251603046473;14=1.027959226511875}
    21 = @{2=1.02;3=1.066409189396018;4=1.070783756911128;5=1.079351749843768;6=1.084495773916378;9=1.052975834411126;10=1.07278332912074;11=1.074239078211258;12=1.082777684772517;13=1.087904311047917;14=1.027495102441007}
    22 = @{2=1.02;3=1.065349329706057;4=1.069930030167643;5=1.078341369455743;6=1.08352115787234;9=1.052659196664748;10=1.071966110520214;11=1.073512613559283;12=1.081893943567106;13=1.087055507395889;14=1.027202082961801}
    23 = @{2=1.02;3=1.06591131889612;4=1.070382742189848;5=1.078877094571059;6=1.084037930625155;9=1.052827271164995;10=1.07239951047043;11=1.073897918472765;12=1.082362586875509;13=1.087505638057563;14=1.027357533150758}
    24 = @{2=1.02;3=1.068120814146057;4=1.072162057802807;5=1.080983941203943;6=1.086070006138785;9=1.053484126746165;10=1.074101876729393;11=1.075410594124336;12=1.084204181906963;13=1.089274162054254;14=1.027966987316889}
    25 = @{2=1.02;3=1.070678873009092;4=1.074220899262881;5=1.083424393205976;6=1.088423349537449;9=1.054236442623044;10=1.076069535056937;11=1.07715738408156;12=1.086334408383341;13=1.091319191136055;14=1.02866904046829}
}

foreach ($r in $data.Keys) {
    foreach ($c in $data[$r].Keys) {
        $ws.Cells.Item($r, $c).Value = $data[$r][$c]
    }
}
